$d = $word.ActiveDocument

$replacements = @(
    @{old="2025-09-20 Saturday"; new="2025-09-21 Sunday"},
    @{old="110÷4="; new="724÷8="},
    @{old="299÷4="; new="356÷4="},
    @{old="632÷2="; new="442÷4="},
    @{old="899÷8="; new="898÷9="},
    @{old="440÷2="; new="596÷3="},
    @{old="709÷8="; new="568÷4="},
    @{old="634÷6="; new="238÷5="},
    @{old="406÷8="; new="273÷6="},
    @{old="806÷4="; new="497÷9="},
    @{old="733÷4="; new="959÷7="},
    @{old="928÷7="; new="387÷5="},
    @{old="348÷5="; new="300÷7="},
    @{old="206÷9="; new="506÷9="},
    @{old="649÷2="; new="610÷4="},
    @{old="689÷8="; new="585÷2="},
    @{old="274÷9="; new="946÷5="},
    @{old="540÷9="; new="746÷9="},
    @{old="513÷9="; new="682÷6="},
    @{old="672÷6="; new="359÷7="},
    @{old="755÷6="; new="768÷4="},
    @{old="481÷4="; new="938÷6="},
    @{old="286÷4="; new="141÷2="},
    @{old="765÷9="; new="774÷6="},
    @{old="817÷3="; new="407÷7="},
    @{old="601÷7="; new="835÷8="}
)

foreach ($r in $replacements) {
    $range = $d.Content
    $range.Find.Execute($r.old, $true, $false, $false, $false, $false, $true, 1, $false, $r.new, 2)
}
